$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2881169905109251
$ws.Range("C2").Value = 9.983522426115931
$ws.Range("D2").Value = 19575605.8673771
$ws.Range("E2").Value = 2459690191846.092
$ws.Range("G2").Value = 2459709767462.231
